# Auto-generated edit script applying cryptos.xlsx diff (18 Jan 2023 symbol-list update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '294.89'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '-2.35%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '31.20'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '-2.37%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '4.924'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-1.66%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.07339'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '-6.39%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.829'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-13.47%'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '7.675'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '-1.79%'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.752'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-0.45%'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9061'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-2.05%'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1653'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '-5.37%'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07620'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-3.17%'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08179'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-6.75%'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.02987'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-4.60%'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09971'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.38%'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001519'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '0.64%'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.005652'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-4.50%'
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.461'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '0.32%'
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'BTSEToken'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.097'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-7.55%'
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'BitpandaEcosystemToken'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.3280'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '0.20%'
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'ProBitToken'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.1306'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '1.13%'
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'MCDex'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.339'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '4.52%'
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'ZBToken'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.2004'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '11.81%'
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'CoinExToken'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04475'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-2.60%'
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'BitKan'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001227'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-1.03%'
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'HotbitToken'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004048'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-9.62%'
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'NitroEx'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0001251'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '0.29%'
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'Spectre.aiUtilityToken'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '--'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '--%'
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'LegolasExchange'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo'
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'BitZToken'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz'
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'Birake'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir'
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'NashExchange'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex'
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'AAXToken'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab'
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'CenX'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V4XJUvLQb+cenx-cenx'
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'BNIXToken'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix'
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'UpBots'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01649'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '-4.67%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04413'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-7.45%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007431'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-0.91%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1322'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002102'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '1.22%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.01105'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '2.84%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005956'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-1.63%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000750'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-0.04%'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.078'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '153.20%'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '-11.62%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002100'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-0.04%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0002000'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.04%'
